$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings: _old -> _FV2410, _new -> _FV2504
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2410")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2504")
}

# Create a table (ListObject) over the data range A1:U64
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split)
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A1").Select()
